$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-detected as a number by
# Excel's type inference: force text format first so the literal string (with
# original formatting, e.g. "577.57") is preserved, matching the source data.
$numericLooking = @("D5", "D6", "D10", "D11", "D12", "D14", "D17", "D20", "D23", "D24", "D25", "D29", "D30", "D31", "D32", "D36", "D37", "D38", "D39", "D40", "D41", "D44", "D45", "D46", "D47", "D48", "D51")
foreach ($ref in $numericLooking) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D2").Value = "61.699.50"
$ws.Range("E2").Value = "  +3.79%  "
$ws.Range("D3").Value = "3.079.63"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "577.57"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("D6").Value = "141.84"
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.067.79"
$ws.Range("E8").Value = "  +2.46%  "
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  +3.67%  "
$ws.Range("D11").Value = "5.47"
$ws.Range("E11").Value = "  +8.01%  "
$ws.Range("D12").Value = "0.466"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("D14").Value = "35.06"
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "3.584.09"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "7.26"
$ws.Range("E17").Value = "  +3.29%  "
$ws.Range("D18").Value = "3.071.24"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("D19").Value = "61.625.16"
$ws.Range("E19").Value = "  +3.88%  "
$ws.Range("D20").Value = "448.65"
$ws.Range("E20").Value = "  +4.18%  "
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("D23").Value = "7.43"
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("D24").Value = "13.55"
$ws.Range("D25").Value = "82.06"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  +4.37%  "
$ws.Range("D29").Value = "2.64"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").Value = "8.06"
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").Value = "6.69"
$ws.Range("E31").Value = "  +8.09%  "
$ws.Range("D32").Value = "26.57"
$ws.Range("E32").Value = "  +3.11%  "
$ws.Range("E33").Value = "  +8.62%  "
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("D35").Value = "0.0₃0793"
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("D36").Value = "6.06"
$ws.Range("E36").Value = "  +4.63%  "
$ws.Range("D37").Value = "2.17"
$ws.Range("E37").Value = "  +3.91%  "
$ws.Range("D38").Value = "50.11"
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("D39").Value = "2.95"
$ws.Range("E39").Value = "  +6.74%  "
$ws.Range("D40").Value = "8.83"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").Value = "425.14"
$ws.Range("E41").Value = "  +3.88%  "
$ws.Range("E42").Value = "  +3.91%  "
$ws.Range("D43").Value = "2.779.78"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.108"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.268"
$ws.Range("E45").Value = "  +6.30%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.10"
$ws.Range("E46").Value = "  +3.53%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "35.30"
$ws.Range("E47").Value = "  +7.95%  "
$ws.Range("D48").Value = "125.07"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").Value = "23.90"
$ws.Range("E51").Value = "  +1.01%  "

# Restore default (Normal) style on those cells so only the value changes,
# not the cell formatting/style index.
foreach ($ref in $numericLooking) { $ws.Range($ref).Style = "Normal" }
